$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Baseline both" column (E) const coefficient and standard
# error, which were previously duplicated from the "Baseline time FE"
# column (D) by mistake.
$ws.Range("E3").Value = "'-0.58***"
$ws.Range("E4").Value = "'(0.02) "

# Move the active selection to I13 (matches the saved sheet view state).
$ws.Range("I13").Select()
